# Refresh crypto price/volume table (GitHub Actions data-refresh commit).
# All data cells in this sheet are stored as text (inline strings), so every
# write below forces NumberFormat "@" on numeric-looking cells (Price /
# Volume(1h) columns) before assigning the value, to stop Excel from
# auto-coercing them into floating point numbers / dates and losing the
# exact textual representation (trailing zeros, percent signs, thousands
# separators, tiny-magnitude decimals, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

function Set-PlainCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Row 2 - BNB
Set-TextCell "D2" "289.97"
Set-TextCell "E2" "-4.12%"

# Row 3 - OKB
Set-TextCell "D3" "30.81"
Set-TextCell "E3" "-4.33%"

# Row 4 - HuobiToken
Set-TextCell "D4" "4.878"
Set-TextCell "E4" "-2.42%"

# Row 5 - Cronos
Set-TextCell "D5" "0.07167"
Set-TextCell "E5" "-9.43%"

# Row 6 - now FTXToken (was KuCoinToken)
Set-PlainCell "B6" "FTXToken"
Set-PlainCell "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell "D6" "1.759"
Set-TextCell "E6" "-16.49%"

# Row 7 - now KuCoinToken (was FTXToken)
Set-PlainCell "B7" "KuCoinToken"
Set-PlainCell "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell "D7" "7.682"
Set-TextCell "E7" "-2.39%"

# Row 8 - GateToken
Set-TextCell "D8" "3.738"
Set-TextCell "E8" "-1.84%"

# Row 9 - MXToken
Set-TextCell "D9" "0.8943"
Set-TextCell "E9" "-3.35%"

# Row 10 - WazirX
Set-TextCell "D10" "0.1667"
Set-TextCell "E10" "-4.84%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-TextCell "D11" "0.07441"
Set-TextCell "E11" "-5.97%"

# Row 12 - MandalaExchangeToken
Set-TextCell "D12" "0.08093"
Set-TextCell "E12" "-7.13%"

# Row 13 - BitrueCoin
Set-TextCell "D13" "0.02975"
Set-TextCell "E13" "-5.11%"

# Row 14 - BitMartToken
Set-TextCell "D14" "0.09996"
Set-TextCell "E14" "-0.34%"

# Row 15 - BitForexToken
Set-TextCell "D15" "0.001494"
Set-TextCell "E15" "-2.00%"

# Row 16 - now TigerCash (was CoinExToken)
Set-PlainCell "B16" "TigerCash"
Set-PlainCell "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D16" "0.005864"
Set-TextCell "E16" "-2.54%"

# Row 17 - now UpBots (was TigerCash)
Set-PlainCell "B17" "UpBots"
Set-PlainCell "C17" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextCell "D17" "0.007527"
Set-TextCell "E17" "2,125.44%"

# Row 18 - now LEO (was UpBots)
Set-PlainCell "B18" "LEO"
Set-PlainCell "C18" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D18" "3.457"
Set-TextCell "E18" "-0.29%"

# Row 19 - now BTSEToken (was LEO)
Set-PlainCell "B19" "BTSEToken"
Set-PlainCell "C19" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell "D19" "2.104"
Set-TextCell "E19" "-7.61%"

# Row 20 - now BitpandaEcosystemToken (was BTSEToken)
Set-PlainCell "B20" "BitpandaEcosystemToken"
Set-PlainCell "C20" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell "D20" "0.3293"
Set-TextCell "E20" "0.18%"

# Row 21 - now ProBitToken (was BitpandaEcosystemToken)
Set-PlainCell "B21" "ProBitToken"
Set-PlainCell "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell "D21" "0.1299"
Set-TextCell "E21" "0.64%"

# Row 22 - now MCDex (was ProBitToken)
Set-PlainCell "B22" "MCDex"
Set-PlainCell "C22" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D22" "4.383"
Set-TextCell "E22" "4.36%"

# Row 23 - now ZBToken (was MCDex)
Set-PlainCell "B23" "ZBToken"
Set-PlainCell "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell "D23" "0.2003"
Set-TextCell "E23" "11.72%"

# Row 24 - now CoinExToken (was ZBToken)
Set-PlainCell "B24" "CoinExToken"
Set-PlainCell "C24" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D24" "0.04470"
Set-TextCell "E24" "-2.88%"

# Row 25 - BitKan
Set-TextCell "D25" "0.001214"
Set-TextCell "E25" "-1.91%"

# Row 26 - HotbitToken
Set-TextCell "D26" "0.004015"
Set-TextCell "E26" "-10.18%"

# Row 27 - NitroEx
Set-TextCell "D27" "0.0001251"
Set-TextCell "E27" "0.10%"

# Row 39 - One
Set-TextCell "D39" "0.01650"
Set-TextCell "E39" "-4.10%"

# Row 40 - IDEX
Set-TextCell "D40" "0.04337"
Set-TextCell "E40" "-9.85%"

# Row 41 - KickToken
Set-TextCell "D41" "0.007413"
Set-TextCell "E41" "-0.08%"

# Row 42 - BKEXToken (price cell unchanged)
Set-TextCell "E42" "-3.86%"

# Row 43 - CEJI
Set-TextCell "D43" "0.002039"
Set-TextCell "E43" "-13.60%"

# Row 44 - LocalTraders
Set-TextCell "D44" "0.01017"
Set-TextCell "E44" "-0.88%"

# Row 45 - CoinLion (price cell unchanged)
Set-TextCell "E45" "-4.57%"

# Row 46 - Kangarootoken
Set-TextCell "D46" "0.00000000751"
Set-TextCell "E46" "0.08%"

# Row 47 - BOLO
Set-TextCell "D47" "2.189"
Set-TextCell "E47" "166.77%"

# Row 48 - CoinbaseStockToken
Set-TextCell "D48" "0.003003"
Set-TextCell "E48" "-11.44%"

# Row 49 - CryptobidCoin
Set-TextCell "D49" "0.00002103"
Set-TextCell "E49" "0.08%"

# Row 50 - SpecialPowerGold
Set-TextCell "D50" "0.0002002"
Set-TextCell "E50" "0.08%"
